$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the "Polkadot" row (row 16) to "Dai"
$ws.Range("A16").Value = "Dai"
$ws.Range("B16").Value = "DAI-USD"

# Update the Market Cap column values (C2:C16) to the refreshed import data
$ws.Range("C2").Value  = 667110329317.8585
$ws.Range("C3").Value  = 215741838217.8916
$ws.Range("C4").Value  = 84537893987.737
$ws.Range("C5").Value  = 34701356782.98159
$ws.Range("C6").Value  = 29596095256.35921
$ws.Range("C7").Value  = 25098292739.90195
$ws.Range("C8").Value  = 13760373926.07715
$ws.Range("C9").Value  = 10211926471.06009
$ws.Range("C10").Value = 10141885921.13766
$ws.Range("C11").Value = 8327370113.840555
$ws.Range("C12").Value = 7169555497.918983
$ws.Range("C13").Value = 6217803225.309694
$ws.Range("C14").Value = 5833348926.360915
$ws.Range("C15").Value = 5592650380.386826
$ws.Range("C16").Value = 5348447881.94701
